$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Result")

$ws.Range("A2").Value = 338
$ws.Range("B2").Value = 250
$ws.Range("C2").Value = 0.01
$ws.Range("D2").Value = 0.028
$ws.Range("E2").Value = 0.02
$ws.Range("F2").Value = 0.007
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0.009
$ws.Range("I2").Value = 0.013
$ws.Range("J2").Value = 150.203
$ws.Range("K2").Value = 170.042
$ws.Range("L2").Value = 251.349
$ws.Range("M2").Value = 222.515
$ws.Range("N2").Value = 148.734
$ws.Range("O2").Value = 130.364

$ws.Range("A3").Value = 160.927
$ws.Range("B3").Value = 250
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0.025
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = -0.013
$ws.Range("G3").Value = -0.022
$ws.Range("H3").Value = -0.007
$ws.Range("I3").Value = -0.011
$ws.Range("J3").Value = 149.638
$ws.Range("K3").Value = 169.528
$ws.Range("L3").Value = 250.577
$ws.Range("M3").Value = 222.088
$ws.Range("N3").Value = 148.406
$ws.Range("O3").Value = 129.841
